$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (file name column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Read the existing data rows (everything below the header row)
$data = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $fname = $ws.Cells.Item($r, 1).Value2
    $val = $ws.Cells.Item($r, 2).Value2

    $letters = ""
    $num = 0
    $loc = ""
    if ($fname -match '^([A-Za-z]+)(\d+)') {
        $letters = $matches[1]
        $num = [int]$matches[2]
        $loc = "$letters$num"
    }

    $data += [PSCustomObject]@{
        FileName = $fname
        Value    = $val
        Letters  = $letters
        Num      = $num
        Loc      = $loc
    }
}

# Sort the electrode locations naturally: by column letter(s), then by row number (A1 .. O15)
$sorted = $data | Sort-Object Letters, Num

# Add the new header for the electrode location column, matching the header style
$ws.Cells.Item(1, 3).Value = "Electrode Locations"
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write the sorted rows back out, filling in column C with the electrode location
for ($i = 0; $i -lt $sorted.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $sorted[$i].FileName
    $ws.Cells.Item($r, 2).Value = $sorted[$i].Value
    $ws.Cells.Item($r, 3).Value = $sorted[$i].Loc
}
